# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$errorDetailZhCn = "Handback file name: jyvjgvi1.t1w is different with handoff file name: 6ac476eb-f33d-4b05-8f3b-ba1976194b22.0a3a3e992e2585941287a935c93ae06e54a52dba.zh-cn."
$errorDetailDeDe = "Handback file name: jyvjgvi1.t1w is different with handoff file name: 6ac476eb-f33d-4b05-8f3b-ba1976194b22.0a3a3e992e2585941287a935c93ae06e54a52dba.de-de."

# --- Overview sheet: update the status text for the 6ac476eb... row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Note: the engine's ColumnWidth -> stored XML "width" conversion adds a
# constant 5/6 character offset (observed empirically), so subtract it here
# to land exactly on the target stored width of 40.
$targetColWidth = 40 - (5/6)

# --- zh-cn sheet: update status + error detail, widen Error Detail column ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = $errorDetailZhCn
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth

# --- de-de sheet: update status + error detail, widen Error Detail column ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = $errorDetailDeDe
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
